# Insert two new rows under "Objectives:" / "Train higher education..." (row 11)
# to add a "Docentes responsáveis:" / "8855158 - Morun Bernardino Neto" entry,
# pushing the rest of the table (Programa resumido: ... Requisitos:) down by 2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12..22 down to 14..24, inheriting row 11's per-column formatting
# (column A = s=1 label style, column B = s=2, column C = s=3) for the two new rows.
$ws.Rows("12:13").Insert()

# New label-only row (column A) -- mirrors the "Avaliação:" / "Requisitos:" pattern.
$ws.Range("A12").Value2 = "Docentes responsáveis:"

# New value row (columns B & C only) -- mirrors the last row's B/C-only pattern.
$ws.Range("B13").Value2 = "8855158 - Morun Bernardino Neto"
$ws.Range("C13").Value2 = "8855158 - Morun Bernardino Neto"

# Drop the unused cells the row-insert auto-filled with inherited (but empty) styles,
# so row 12 only has A12 and row 13 only has B13/C13 -- matching the source layout.
$ws.Range("B12:C12").Clear()
$ws.Range("A13").Clear()
